$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the combined order list in its new canonical order with updated
# (combined) quantities, as produced by the CLI's combine_orders routine.
$ws.Range("A1").Value = "Copper Horse - Warhorse 5lb"
$ws.Range("B1").Value = 5

$ws.Range("A2").Value = "Copper Horse - Warhorse Blend (12oz)"
$ws.Range("B2").Value = 8

$ws.Range("A3").Value = "Copper Horse - Rumble Pony (12oz)"
$ws.Range("B3").Value = 7

$ws.Range("A4").Value = "Copper Horse - Clocktower Espresso (12oz)"
$ws.Range("B4").Value = 7

$ws.Range("A5").Value = "Copper Horse - Carriage House Blend (12oz)"
$ws.Range("B5").Value = 13
